$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Total" row correct/total marks
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
